$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 3.2
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("AT2").Value = 2.63
$ws.Range("AX2").Value = 17
$ws.Range("BB2").Value = 201
